# "Generate Report for Archive"
# 1) Status text "Ready for handoff" -> "In Translation" (all sheets/columns that show it)
# 2) Narrow the "status" columns (Overview!E:F, zh-cn!C, de-de!C) from ~17.22 chars to ~13.41 chars

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status values ---
$wsOverview.Range("E2:E4").Value2 = "In Translation"
$wsOverview.Range("F2:F4").Value2 = "In Translation"
$wsZhCn.Range("C2:C4").Value2 = "In Translation"
$wsDeDe.Range("C2:C4").Value2 = "In Translation"

# --- Resize status columns ---
# Target stored OOXML column width is ~13.4101845877511 characters. The
# COM ColumnWidth setter snaps to a whole-pixel grid, so 12.5 is the
# closest settable value that lands nearest the target width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
